$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 100 (pushes old rows 100-200 down to 102-202)
$ws.Rows("100:101").Insert()

# New row 100: Camote, Primera, caja 18 kilos
$ws.Range("A100").Value = 9
$ws.Range("B100").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C100").Value = "Metropolitana"
$ws.Range("D100").Value = 45264
$ws.Range("D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E100").Value = 13
$ws.Range("F100").Value = 100114002
$ws.Range("G100").Value = "Camote"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 970
$ws.Range("K100").Value = 13000
$ws.Range("L100").Value = 14000
$ws.Range("M100").Value = 13485
$ws.Range("N100").Value = "$/caja 18 kilos"
$ws.Range("O100").Value = "Perú"
$ws.Range("P100").Value = 749
$ws.Range("Q100").Value = 18
$ws.Range("R100").Value = "Hortaliza"

# New row 101: Camote, Primera, malla 18 kilos
$ws.Range("A101").Value = 9
$ws.Range("B101").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C101").Value = "Metropolitana"
$ws.Range("D101").Value = 45264
$ws.Range("D101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E101").Value = 13
$ws.Range("F101").Value = 100114002
$ws.Range("G101").Value = "Camote"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 700
$ws.Range("K101").Value = 10000
$ws.Range("L101").Value = 11000
$ws.Range("M101").Value = 10500
$ws.Range("N101").Value = "$/malla 18 kilos"
$ws.Range("O101").Value = "Perú"
$ws.Range("P101").Value = 583
$ws.Range("Q101").Value = 18
$ws.Range("R101").Value = "Hortaliza"
